# Atualização automática: 2025-08-30 09:01:04
#
# This script applies the row-level data changes described by the diff
# against the "dashboard_data" worksheet (sheet1). Rows 7-11 contain a
# cyclical re-ordering of the same records (row 11 -> row 7, and rows
# 7,8,9,10 shift down to 8,9,10,11), and rows 16-17 get corrected
# detection-image / bounding-box values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "First_Coords" (I) and "First_Confidence" (J) columns store
# values that look numeric (e.g. "0.76", or "702,633,740,690" which
# Excel would otherwise read as a comma-thousands number) but must
# stay text, matching the rest of the workbook. Mark the destination
# cells as Text ahead of time so Excel keeps them as strings instead
# of re-interpreting them as numbers when we assign the new values.
$ws.Range("I7:I11").NumberFormat = "@"
$ws.Range("J7:J11").NumberFormat = "@"
$ws.Range("I16:I17").NumberFormat = "@"

# ---- Row 7 (becomes the record that used to be on row 11) ----
$ws.Range("A7").Value2 = "2117575c-4ae1-458c-b88a-fc40f40debdb"
$ws.Range("D7").Value2 = "image_20250727074723_ppp0.jpg"
$ws.Range("E7").Value2 = "PLACA_20250723145134"
$ws.Range("F7").Value2 = "Moura"
$ws.Range("G7").Value2 = 38.06587
$ws.Range("H7").Value2 = -7.221796
$ws.Range("I7").Value2 = "1490,161,1563,258"
$ws.Range("J7").Value2 = "0.62"

# ---- Row 8 (becomes the record that used to be on row 7) ----
$ws.Range("A8").Value2 = "283b6eda-9c83-4cdd-9524-c7c394f2dc89"
$ws.Range("I8").Value2 = "962,713,1006,765"
$ws.Range("J8").Value2 = "0.76"

# ---- Row 9 (becomes the record that used to be on row 8) ----
$ws.Range("A9").Value2 = "a19b65d1-6f97-4841-9e1c-7446a9be92b6"
$ws.Range("I9").Value2 = "967,614,1002,659"
$ws.Range("J9").Value2 = "0.73"

# ---- Row 10 (becomes the record that used to be on row 9) ----
$ws.Range("A10").Value2 = "4be1b1cf-d480-453e-b5fb-d4ecd6764c4d"
$ws.Range("I10").Value2 = "702,633,740,690"
$ws.Range("J10").Value2 = "0.72"

# ---- Row 11 (becomes the record that used to be on row 10) ----
$ws.Range("A11").Value2 = "dfd476d4-7689-4671-a076-78fe3ce806bb"
$ws.Range("D11").Value2 = "image_20250728214139_ppp0.jpg"
$ws.Range("E11").Value2 = "PLACA_20250717165933"
$ws.Range("F11").Value2 = "Beja"
$ws.Range("G11").Value2 = 38.02035
$ws.Range("H11").Value2 = -7.94715
$ws.Range("I11").Value2 = "1254,850,1294,895"
$ws.Range("J11").Value2 = "0.67"

# ---- Row 16: corrected detection image + bounding box ----
$ws.Range("D16").Value2 = "image_20250807111026_ppp0.jpg"
$ws.Range("I16").Value2 = "641,529,688,576"

# ---- Row 17: corrected detection image + bounding box ----
$ws.Range("D17").Value2 = "image_20250807111026_ppp0.jpg"
$ws.Range("I17").Value2 = "793,481,831,526"
